$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 20-21; existing rows 20-27 shift down to 22-29.
$ws.Rows("20:21").Insert()

# Row 20: new Damasco / Castle Brite / Primera entry (Provincia de San Felipe de Aconcagua)
$ws.Range("A20").Value = 10
$ws.Range("B20").Value = "Vega Modelo de Temuco"
$ws.Range("C20").Value = "La Araucanía"
$ws.Range("D20").Value = 44529
$ws.Range("E20").Value = 9
$ws.Range("F20").Value = "Fruta"
$ws.Range("G20").Value = 100103
$ws.Range("H20").Value = "Frutos de hueso (carozo)"
$ws.Range("I20").Value = 100103003
$ws.Range("J20").Value = "Damasco"
$ws.Range("K20").Value = "Castle Brite"
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 75
$ws.Range("N20").Value = 20000
$ws.Range("O20").Value = 20000
$ws.Range("P20").Value = 20000
$ws.Range("Q20").Value = "$/bandeja 10 kilos"
$ws.Range("R20").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S20").Value = 2000
$ws.Range("T20").Value = 10

# Row 21: new Damasco / Castle Brite / Segunda entry (Provincia de San Felipe de Aconcagua)
$ws.Range("A21").Value = 10
$ws.Range("B21").Value = "Vega Modelo de Temuco"
$ws.Range("C21").Value = "La Araucanía"
$ws.Range("D21").Value = 44529
$ws.Range("E21").Value = 9
$ws.Range("F21").Value = "Fruta"
$ws.Range("G21").Value = 100103
$ws.Range("H21").Value = "Frutos de hueso (carozo)"
$ws.Range("I21").Value = 100103003
$ws.Range("J21").Value = "Damasco"
$ws.Range("K21").Value = "Castle Brite"
$ws.Range("L21").Value = "Segunda"
$ws.Range("M21").Value = 45
$ws.Range("N21").Value = 15000
$ws.Range("O21").Value = 15000
$ws.Range("P21").Value = 15000
$ws.Range("Q21").Value = "$/bandeja 10 kilos"
$ws.Range("R21").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S21").Value = 1500
$ws.Range("T21").Value = 10
